# Delegation of assignments for E2 and update of my diary.
#
# Appends a "Day3:" diary entry after the existing text that ends with
# "...Iteration(E2) tomorrow." The new content is:
#   - one blank paragraph
#   - a paragraph with the Day3 entry (several line-broken sentences,
#     a misspelled "sequens" flagged via proofErr, and a trailing
#     manual line break) that also now carries the relocated _GoBack
#     bookmark.

$d = $word.ActiveDocument

# The _GoBack bookmark currently sits at the end of the last paragraph;
# it needs to move into the new paragraph we are about to create, so
# drop it here and re-insert it (via raw XML) at its new home below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$rng = $d.Content
$rng.Collapse(0)

$w = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'
$lang = "<w:rPr><w:lang w:val='en-US'/></w:rPr>"

$frag  = "<w:p xmlns:w='$w'><w:pPr>$lang</w:pPr></w:p>"
$frag += "<w:p xmlns:w='$w'><w:pPr>$lang</w:pPr>"
$frag += "<w:r>$lang<w:t>Day3:</w:t></w:r>"
$frag += "<w:r>$lang<w:br/><w:t>Start of E2.</w:t></w:r>"
$frag += "<w:r>$lang<w:br/><w:t>The programming assignments are far more difficult to distribute, than the last iteration.</w:t></w:r>"
$frag += "<w:r>$lang<w:br/><w:t xml:space='preserve'>So instead, I have distributed the assignments of making the Fully dressed use cases and </w:t></w:r>"
$frag += "<w:proofErr w:type='spellStart'/>"
$frag += "<w:r>$lang<w:t>sequens</w:t></w:r>"
$frag += "<w:proofErr w:type='spellEnd'/>"
$frag += "<w:r>$lang<w:t xml:space='preserve'> diagrams.</w:t></w:r>"
$frag += "<w:r>$lang<w:br/><w:t>The plan is that these are finalized tomorrow. Then we can talk about how we distribute the programming part.</w:t></w:r>"
$frag += "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>"
$frag += "<w:r>$lang<w:br/></w:r>"
$frag += "</w:p>"

$rng.InsertXML($frag)

$d.Save()
